# Apply corrected error-estimation / projected-years results
$wb = $excel.ActiveWorkbook

# --- Sheet "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("C2").Value = 1
$ws1.Range("E2").Value = 3.4

$ws1.Range("C3").Value = 5
$ws1.Range("E3").Value = 17.2

$ws1.Range("B4").Value = 4
$ws1.Range("C4").Value = 19
$ws1.Range("D4").Value = 44.4
$ws1.Range("E4").Value = 65.5

$ws1.Range("C5").Value = 2
$ws1.Range("D5").Value = 44.4
$ws1.Range("E5").Value = 6.9

$ws1.Range("D6").Value = 11.1
$ws1.Range("E6").Value = 6.9

$ws1.Range("B7").Value = 37
$ws1.Range("C7").Value = 125

# --- Sheet "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("C3").Value = 9
$ws4.Range("C4").Value = 29
